$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2025-12-05 Friday" "2025-12-06 Saturday"

Replace-Text "580÷3=193, 1" "521÷3=173, 2"
Replace-Text "384÷2=192, 0" "318÷9=35, 3"
Replace-Text "960÷9=106, 6" "513÷3=171, 0"
Replace-Text "365÷8=45, 5" "966÷5=193, 1"
Replace-Text "767÷9=85, 2" "765÷2=382, 1"

Replace-Text "926÷6=154, 2" "420÷9=46, 6"
Replace-Text "770÷4=192, 2" "540÷3=180, 0"
Replace-Text "467÷6=77, 5" "241÷6=40, 1"
Replace-Text "962÷2=481, 0" "465÷4=116, 1"
Replace-Text "395÷6=65, 5" "104÷4=26, 0"

Replace-Text "574÷8=71, 6" "187÷7=26, 5"
Replace-Text "594÷6=99, 0" "556÷2=278, 0"
Replace-Text "297÷9=33, 0" "826÷3=275, 1"
Replace-Text "897÷4=224, 1" "386÷3=128, 2"
Replace-Text "164÷5=32, 4" "166÷6=27, 4"

Replace-Text "466÷2=233, 0" "156÷6=26, 0"
Replace-Text "671÷8=83, 7" "392÷2=196, 0"
Replace-Text "225÷8=28, 1" "564÷5=112, 4"
Replace-Text "360÷8=45, 0" "632÷2=316, 0"
Replace-Text "387÷4=96, 3" "202÷7=28, 6"

Replace-Text "356÷4=89, 0" "707÷6=117, 5"
Replace-Text "573÷6=95, 3" "822÷3=274, 0"
Replace-Text "278÷5=55, 3" "596÷8=74, 4"
Replace-Text "192÷8=24, 0" "553÷9=61, 4"
Replace-Text "866÷5=173, 1" "280÷5=56, 0"

Write-Output "done"
